$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.322
$ws.Range("D7").Value = -7.699
$ws.Range("B10").Value = 5.873
$ws.Range("B12").Value = 5.089999999999999
$ws.Range("D15").Value = -8.175000000000001
$ws.Range("B18").Value = 5.266
$ws.Range("E18").Value = 16.551
$ws.Range("E19").Value = 16.538
$ws.Range("D20").Value = -7.57
$ws.Range("E27").Value = 16.257
$ws.Range("D29").Value = -7.282000000000001
$ws.Range("D30").Value = -7.236
$ws.Range("D31").Value = -8.266999999999999
$ws.Range("B37").Value = 8.73
$ws.Range("D40").Value = -7.93
$ws.Range("E42").Value = 16.555
$ws.Range("E44").Value = 16.623
$ws.Range("E47").Value = 16.467
$ws.Range("B55").Value = 4.693
$ws.Range("E58").Value = 16.387
$ws.Range("B68").Value = 5.598000000000001
$ws.Range("D68").Value = -6.872
$ws.Range("E73").Value = 16.635
$ws.Range("D76").Value = -7.753000000000002
$ws.Range("B77").Value = 6.123
$ws.Range("B78").Value = 7.811
$ws.Range("D87").Value = -8.113
$ws.Range("D88").Value = -7.916000000000001
$ws.Range("E95").Value = 17.499
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.228
$ws.Range("D101").Value = -7.9
$ws.Range("E101").Value = 16.511
$ws.Range("D102").Value = -8.036
